$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-fill the Cost / Unit Cost / Kit Cost / Sales Tax formula columns ---
# so they become Excel "shared formula" groups spanning the whole table
# (rows 3:25), matching how a fill-down over the existing data looks.
$ws.Range("G3:G25").Formula = "=B3+D3+F3"
$ws.Range("H3:H25").Formula = "=G3/C3"
$ws.Range("J3:J25").Formula = "=H3*I3"
$ws.Range("F4:F25").Formula = "=(B4+D4+E4)*0.1075"

# --- Row 13: "2.1mmx5.5mm DC barrel jack" price/qty change ---
# Purchase price drops from 9.99 to 0.88, purchase qty from 10 to 1,
# and a new Tariff (column E) formula is introduced for this row.
$ws.Range("B13").Value = 0.88
$ws.Range("C13").Value = 1
$ws.Range("E13").Formula = "=B13*0.1"

# --- Insert a new row 26 for the new "audio jack" part ---
# This pushes the old blank spacer row (26) down to 27 and the
# Totals row (27) down to 28; Excel auto-updates the SUM() ranges.
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "audio jack"
$ws.Range("B26").Value = 11.69
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 0
$ws.Range("E26").Formula = "=B25*0.1"
$ws.Range("F26").Formula = "=(B26+D26+E26)*0.1075"
$ws.Range("G26").Formula = "=B26+D26+F26"
$ws.Range("H26").Formula = "=G26/C26"
$ws.Range("I26").Value = 1
$ws.Range("J26").Formula = "=H26*I26"
$ws.Range("K26").Value = "https://www.digikey.com/en/products/detail/cui-devices/SJ1-3533NG/738701"

# --- Restore the UI selection state to match the post-edit workbook ---
[void]$ws.Range("E31").Select()
